$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (want-to-go count) values
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F5").Value = 5135
$wsExhibition.Range("F14").Value = 4084
$wsExhibition.Range("F16").Value = 157
$wsExhibition.Range("F17").Value = 142
$wsExhibition.Range("F19").Value = 3083
$wsExhibition.Range("F23").Value = 41
$wsExhibition.Range("F25").Value = 84

# Sheet "全部类型" (All types) - same events, shifted by one row - update the same counts
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value = 5135
$wsAll.Range("F15").Value = 4084
$wsAll.Range("F17").Value = 157
$wsAll.Range("F18").Value = 142
$wsAll.Range("F20").Value = 3083
$wsAll.Range("F24").Value = 41
$wsAll.Range("F26").Value = 84
